# Update sales.xlsx from Streamlit app
# - B3's BARCODE value was saved as text; correct it to a real number (5421)
# - Append a new sale row (row 4) written by the app

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 3: BARCODE should be numeric, not text ---
$ws.Range("B3").Value = 5421

# --- Append row 4 ---
$ws.Range("A4").Value = "2025-09-11 06:08:50"

# BARCODE (B4) and Product (C4) must stay text even though they look numeric,
# so force a text number format before assigning, then drop the formatting
# again so the cell keeps default styling but retains its text data type.
$cellB4 = $ws.Range("B4")
$cellB4.NumberFormat = "@"
$cellB4.Value = "1220"
$cellB4.ClearFormats()

$cellC4 = $ws.Range("C4")
$cellC4.NumberFormat = "@"
$cellC4.Value = "6698"
$cellC4.ClearFormats()

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 199
$ws.Range("F4").Value = "ALLAN"
$ws.Range("G4").Value = "CHRIS"
$ws.Range("H4").Value = "Sale"
